{"js": "// Design Doc.docx edit:\n//  1. \"ORDER:\" paragraph \u2014 drop the spell-check proofErr wrapper around the run.\n//  2. \"Block health display real angle vs fast angle\" paragraph's text is replaced\n//     with \"Display settings separate from gameplay settings in load files\",\n//     and the four paragraphs that followed it (\"Block health display health left\n//     to right\", \"Option to not display any of them\", \"Block health percentage and\n//     option\", \"Can have any or none of the health displays\") are removed.\n//  3. \"UIState and UIPanel ShowHideUI\" paragraph \u2014 collapse the multiple runs\n//     (separated by spell-check proofErr wrappers) into a single plain run.\n\nconst body = context.document.body;\n\n// --- 1. \"ORDER:\" \u2014 remove proofErr spellStart/spellEnd wrapper -------------\nconst orderResults = body.search(\"ORDER:\", { matchCase: true, matchWholeWord: true });\norderResults.load(\"items\");\nawait context.sync();\n\nif (orderResults.items.length > 0) {\n  const orderPara = orderResults.items[0].paragraphs.getFirst();\n  // Re-inserting a clean paragraph (then dropping the old one) rebuilds the\n  // run without the spell-check proofErr markers that wrapped it before.\n  orderPara.insertParagraph(\"ORDER:\", \"Before\");\n  orderPara.delete();\n  await context.sync();\n}\n\n// --- 2. Replace text + delete the four paragraphs that followed ------------\nconst angleResults = body.search(\"Block health display real angle vs fast angle\", { matchCase: true });\nangleResults.load(\"items\");\nawait context.sync();\n\nif (angleResults.items.length > 0) {\n  const targetPara = angleResults.items[0].paragraphs.getFirst();\n  targetPara\n    .getRange()\n    .insertText(\"Display settings separate from gameplay settings in load files\", \"Replace\");\n  await context.sync();\n}\n\nconst paragraphsToRemove = [\n  \"Block health display health left to right\",\n  \"Option to not display any of them\",\n  \"Block health percentage and option\",\n  \"Can have any or none of the health displays\",\n];\n\nfor (const text of paragraphsToRemove) {\n  const results = body.search(text, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].paragraphs.getFirst().delete();\n    await context.sync();\n  }\n}\n\n// --- 3. \"UIState and UIPanel ShowHideUI\" \u2014 collapse runs, drop proofErr ----\nconst uiResults = body.search(\"UIState and UIPanel ShowHideUI\", { matchCase: true });\nuiResults.load(\"items\");\nawait context.sync();\n\nif (uiResults.items.length > 0) {\n  const uiPara = uiResults.items[0].paragraphs.getFirst();\n  uiPara.insertParagraph(\"UIState and UIPanel ShowHideUI\", \"Before\");\n  uiPara.delete();\n  await context.sync();\n}\n", "ps1": "# Design Doc.docx edit:\n#  1. \"ORDER:\" paragraph \u2014 drop the spell-check proofErr wrapper around the run.\n#  2. \"Block health display real angle vs fast angle\" paragraph's text is replaced\n#     with \"Display settings separate from gameplay settings in load files\",\n#     and the four paragraphs that followed it (\"Block health display health left\n#     to right\", \"Option to not display any of them\", \"Block health percentage and\n#     option\", \"Can have any or none of the health displays\") are removed.\n#  3. \"UIState and UIPanel ShowHideUI\" paragraph \u2014 collapse the multiple runs\n#     (separated by spell-check proofErr wrappers) into a single plain run.\n\n$d = $word.ActiveDocument\n\n# Re-inserting a clean paragraph (then dropping the old one, mark included)\n# rebuilds the run without any spell-check proofErr markers that wrapped it\n# before -- a plain Range.Text assignment leaves proofErr markers in place.\nfunction Replace-ParagraphCleanly($searchText, $newText) {\n  $r = $d.Content\n  $found = $r.Find.Execute($searchText)\n  if (-not $found) { return }\n  $r.Expand(4)  # wdParagraph -- include the paragraph mark\n  $oldStart = $r.Start\n  $oldEnd = $r.End\n  $insertStr = $newText + \"`r\"\n  $r.InsertBefore($insertStr)\n  $oldRange = $d.Range($oldStart + $insertStr.Length, $oldEnd + $insertStr.Length)\n  $oldRange.Delete()\n}\n\nfunction Delete-ParagraphWithText($searchText) {\n  $r = $d.Content\n  $found = $r.Find.Execute($searchText)\n  if (-not $found) { return }\n  $r.Expand(4)  # wdParagraph -- include the paragraph mark\n  $r.Delete()\n}\n\n# --- 1. \"ORDER:\" \u2014 remove proofErr spellStart/spellEnd wrapper -------------\nReplace-ParagraphCleanly \"ORDER:\" \"ORDER:\"\n\n# --- 2. Replace text + delete the four paragraphs that followed ------------\n$r = $d.Content\n$found = $r.Find.Execute(\"Block health display real angle vs fast angle\")\nif ($found) {\n  $r.Text = \"Display settings separate from gameplay settings in load files\"\n}\n\nDelete-ParagraphWithText \"Block health display health left to right\"\nDelete-ParagraphWithText \"Option to not display any of them\"\nDelete-ParagraphWithText \"Block health percentage and option\"\nDelete-ParagraphWithText \"Can have any or none of the health displays\"\n\n# --- 3. \"UIState and UIPanel ShowHideUI\" \u2014 collapse runs, drop proofErr ----\nReplace-ParagraphCleanly \"UIState and UIPanel ShowHideUI\" \"UIState and UIPanel ShowHideUI\"\n"}
